$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S2").Value = "R`$ 165.83"
$ws.Range("U2").Value = "27/04/2023 15:23:27"
$ws.Range("D3").Value = "Pedacinho do Céu"
$ws.Range("I3").Value = "1 cama de casal"
$ws.Range("K3").Value = "R`$85 por noite, originalmente R`$93"
$ws.Range("O3").Value = "4,92 (12)"
$ws.Range("S3").Value = "R`$ 165.83"
$ws.Range("U3").Value = "27/04/2023 15:23:27"
$ws.Range("A4").Value = "Casa em Parque Hotel"
$ws.Range("D4").Value = "Aconchegante Casa - Quintal & Garagem"
$ws.Range("I4").Value = "1 cama queen"
$ws.Range("K4").Value = "R`$137 por noite"
$ws.Range("O4").Value = "4,89 (46)"
$ws.Range("S4").Value = "R`$ 165.83"
$ws.Range("U4").Value = "27/04/2023 15:23:27"
$ws.Range("A5").Value = "Casa em Fazendinha"
$ws.Range("D5").Value = "Casa em Araruama 2"
$ws.Range("I5").Value = "4 camas"
$ws.Range("K5").Value = "R`$98 por noite, originalmente R`$109"
$ws.Range("O5").Value = "5,0 (27)"
$ws.Range("Q5").Value = "Superhost"
$ws.Range("S5").Value = "R`$ 165.83"
$ws.Range("U5").Value = "27/04/2023 15:23:27"
$ws.Range("D6").Value = "Casa em Araruama 1"
$ws.Range("I6").Value = "4 camas"
$ws.Range("K6").Value = "R`$106 por noite"
$ws.Range("O6").Value = "4,97 (159)"
$ws.Range("Q6").Value = "Superhost"
$ws.Range("S6").Value = "R`$ 165.83"
$ws.Range("U6").Value = "27/04/2023 15:23:27"
$ws.Range("A7").Value = "Casa em Araruama"
$ws.Range("D7").Value = "Meu aconchego"
$ws.Range("I7").Value = "9 camas"
$ws.Range("K7").Value = "R`$202 por noite, originalmente R`$222"
$ws.Range("O7").Value = "4,9 (10)"
$ws.Range("Q7").Value = $null
$ws.Range("S7").Value = "R`$ 165.83"
$ws.Range("U7").Value = "27/04/2023 15:23:27"
$ws.Range("A8").Value = "Apartamento em Araruama"
$ws.Range("D8").Value = "Apartamento em frente a lagoa de Araruama"
$ws.Range("I8").Value = "2 camas"
$ws.Range("K8").Value = "R`$173 por noite, originalmente R`$211"
$ws.Range("O8").Value = "4,85 (26)"
$ws.Range("Q8").Value = $null
$ws.Range("S8").Value = "R`$ 165.83"
$ws.Range("U8").Value = "27/04/2023 15:23:27"
$ws.Range("S9").Value = "R`$ 165.83"
$ws.Range("U9").Value = "27/04/2023 15:23:27"
$ws.Range("A10").Value = "Casa em Araruama"
$ws.Range("D10").Value = "loft <SPA< onde você descansa relaxa e se renova"
$ws.Range("I10").Value = "2 camas"
$ws.Range("K10").Value = "R`$189 por noite, originalmente R`$302"
$ws.Range("O10").Value = "5,0 (9)"
$ws.Range("Q10").Value = "Superhost"
$ws.Range("S10").Value = "R`$ 165.83"
$ws.Range("U10").Value = "27/04/2023 15:23:27"
$ws.Range("D11").Value = "Loft completo para temporada"
$ws.Range("I11").Value = "2 camas"
$ws.Range("K11").Value = "R`$296 por noite"
$ws.Range("O11").Value = "4,92 (25)"
$ws.Range("Q11").Value = $null
$ws.Range("S11").Value = "R`$ 165.83"
$ws.Range("U11").Value = "27/04/2023 15:23:27"
$ws.Range("A12").Value = "Casa em Araruama"
$ws.Range("D12").Value = "Aluguel por temporada com piscina privativa"
$ws.Range("I12").Value = "3 camas"
$ws.Range("K12").Value = "R`$315 por noite"
$ws.Range("O12").Value = "4,97 (67)"
$ws.Range("Q12").Value = "Superhost"
$ws.Range("S12").Value = "R`$ 165.83"
$ws.Range("U12").Value = "27/04/2023 15:23:27"
$ws.Range("A13").Value = "Apartamento em Araruama"
$ws.Range("D13").Value = "Apto na Região dos lagos. Aconchegante e central."
$ws.Range("I13").Value = "2 camas"
$ws.Range("K13").Value = "R`$140 por noite"
$ws.Range("O13").Value = "4,71 (7)"
$ws.Range("S13").Value = "R`$ 165.83"
$ws.Range("U13").Value = "27/04/2023 15:23:27"
$ws.Range("A14").Value = "Casa de campo em Praia Seca"
$ws.Range("D14").Value = "Casa Maravilhosa com Lagoa privativa"
$ws.Range("I14").Value = "3 camas de casal"
$ws.Range("K14").Value = "R`$117 por noite"
$ws.Range("O14").Value = "4,77 (22)"
$ws.Range("Q14").Value = $null
$ws.Range("S14").Value = "R`$ 165.83"
$ws.Range("U14").Value = "27/04/2023 15:23:27"
$ws.Range("S15").Value = "R`$ 165.83"
$ws.Range("U15").Value = "27/04/2023 15:23:27"
$ws.Range("A16").Value = "Casa de hóspedes em Pontinha"
$ws.Range("D16").Value = "Casa do Alto da Pontinha"
$ws.Range("K16").Value = "R`$72 por noite"
$ws.Range("O16").Value = "4,83 (18)"
$ws.Range("S16").Value = "R`$ 165.83"
$ws.Range("U16").Value = "27/04/2023 15:23:27"
$ws.Range("A17").Value = "Casa de campo em Outeiro"
$ws.Range("D17").Value = "Casa com piscina no Condomínio Sonho de Vida"
$ws.Range("I17").Value = "4 camas"
$ws.Range("K17").Value = "R`$290 por noite, originalmente R`$352"
$ws.Range("O17").Value = "5,0 (22)"
$ws.Range("Q17").Value = "Superhost"
$ws.Range("S17").Value = "R`$ 165.83"
$ws.Range("U17").Value = "27/04/2023 15:23:27"
$ws.Range("A18").Value = "Apartamento em Araruama"
$ws.Range("D18").Value = "Kitinete em Araruama/RJ!`n`nCentro!`n`nPerto da Praia!"
$ws.Range("I18").Value = "2 camas"
$ws.Range("K18").Value = "R`$198 por noite"
$ws.Range("O18").Value = "4,8 (10)"
$ws.Range("Q18").Value = $null
$ws.Range("S18").Value = "R`$ 165.83"
$ws.Range("U18").Value = "27/04/2023 15:23:27"
$ws.Range("A19").Value = "Apartamento em Parque Hotel"
$ws.Range("D19").Value = "Apartamento na Lagoa de Araruama"
$ws.Range("I19").Value = "3 camas"
$ws.Range("K19").Value = "R`$208 por noite"
$ws.Range("O19").Value = "4,94 (17)"
$ws.Range("Q19").Value = "Superhost"
$ws.Range("S19").Value = "R`$ 165.83"
$ws.Range("U19").Value = "27/04/2023 15:23:27"
